# Updates to latest 4.0
$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")

# Update the conversion factor (About!A26) and its accompanying label (About!B26)
$about.Range("A26").Value = 0.75350342301658668
$about.Range("B26").Value = "2023 dollars per 2012 dollar"

# Move the active selection to A26 to match the saved view state
$about.Activate()
$about.Range("A26").Select()

$wb.Save()
